$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first worker's block (ANGELICA MARIA RODRIGUEZ TORRES, rows 16-21).
# This shifts the second worker's (CARMEN ROCIO BALLESTEROS FLOREZ) 8 rows of
# period data up to rows 16-23, and the footer block up from rows 34-35 to 28-29,
# dropping the now-empty trailing rows.
$ws.Rows("16:21").Delete()

# Update the summary figures for the new single-worker statement.
$ws.Range("E11").Value = 254623
$ws.Range("C13").Value = 1

# Reorder CARMEN's period/mora-value rows into ascending period order
# (1909 .. 2006) as part of the "parte 1" rebuild of the statement data.
$ws.Range("E16").Value = "1909"
$ws.Range("F16").Value = 31249

$ws.Range("E17").Value = "1910"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "2001"
$ws.Range("F18").Value = 33125

$ws.Range("E19").Value = "2002"
$ws.Range("F19").Value = 33125

$ws.Range("E20").Value = "2003"
$ws.Range("F20").Value = 33125

$ws.Range("E21").Value = "2004"
$ws.Range("F21").Value = 33125

$ws.Range("E22").Value = "2005"
$ws.Range("F22").Value = 33125

$ws.Range("E23").Value = "2006"
$ws.Range("F23").Value = 26500
